$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 1118865.4
$ws.Range("I34").Value = 1252473.5
$ws.Range("K34").Value = 1252473.5
$ws.Range("M34").Value = -1252270.5
$ws.Range("H36").Value = 1118865.4
$ws.Range("I36").Value = 1252473.5
$ws.Range("K36").Value = 1252473.5
$ws.Range("M36").Value = -1251758.5
$ws.Range("H70").Value = 2543.7896
$ws.Range("I70").Value = 1798.3334
$ws.Range("J70").Value = 2683.5625
$ws.Range("K70").Value = 5395.0002
$ws.Range("L70").Value = 8050.6875
$ws.Range("M70").Value = -5125.0002
$ws.Range("N70").Value = -8590.6875
$ws.Range("H73").Value = 2543.7896
$ws.Range("I73").Value = 1798.3334
$ws.Range("J73").Value = 2683.5625
$ws.Range("K73").Value = 5395.0002
$ws.Range("L73").Value = 8050.6875
$ws.Range("M73").Value = -4459.0002
$ws.Range("N73").Value = -9922.6875
$ws.Range("H116").Value = 3107.8333
$ws.Range("I116").Value = 3599.3333
$ws.Range("J116").Value = 1633.3334
$ws.Range("K116").Value = 3599.3333
$ws.Range("L116").Value = 1633.3334
$ws.Range("M116").Value = -157.3332999999998
$ws.Range("N116").Value = -8517.3334
$ws.Range("H129").Value = 1195.4375
$ws.Range("J129").Value = 1255.3334
$ws.Range("L129").Value = 3766.0002
$ws.Range("N129").Value = -13766.0002
$ws.Range("H133").Value = 49990
$ws.Range("J133").Value = 49990
$ws.Range("L133").Value = 49990
$ws.Range("N133").Value = -60110
$ws.Range("H138").Value = 8407252
$ws.Range("I138").Value = 3574966.8
$ws.Range("J138").Value = 10420704
$ws.Range("K138").Value = 10724900.4
$ws.Range("L138").Value = 31262112
$ws.Range("M138").Value = -10719760.4
$ws.Range("N138").Value = -31272392

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 506255.5
$ws.Range("I2").Value = 506255.5
$ws.Range("K2").Value = 506255.5
$ws.Range("M2").Value = -506142.5
$ws.Range("H45").Value = 1346.5
$ws.Range("I45").Value = 1353.1052
$ws.Range("K45").Value = 1353.1052
$ws.Range("M45").Value = -976.1052
$ws.Range("H61").Value = 2921.6785
$ws.Range("I61").Value = 2270.5
$ws.Range("J61").Value = 3789.9167
$ws.Range("K61").Value = 2270.5
$ws.Range("L61").Value = 3789.9167
$ws.Range("M61").Value = -2058.5
$ws.Range("N61").Value = -4213.9167
$ws.Range("H74").Value = 6051.148
$ws.Range("I74").Value = 1114.0526
$ws.Range("K74").Value = 1114.0526
$ws.Range("M74").Value = -240.0526
$ws.Range("H77").Value = 6051.148
$ws.Range("I77").Value = 1114.0526
$ws.Range("K77").Value = 5570.263
$ws.Range("M77").Value = -1202.263
$ws.Range("H97").Value = 6527.9414
$ws.Range("I97").Value = 7676.7856
$ws.Range("K97").Value = 7676.7856
$ws.Range("M97").Value = -7180.7856
$ws.Range("H116").Value = 506255.5
$ws.Range("I116").Value = 506255.5
$ws.Range("K116").Value = 506255.5
$ws.Range("M116").Value = -503961.5
$ws.Range("H133").Value = 48419
$ws.Range("J133").Value = 48419
$ws.Range("L133").Value = 48419
$ws.Range("N133").Value = -53479
$ws.Range("H136").Value = 2921.6785
$ws.Range("I136").Value = 2270.5
$ws.Range("J136").Value = 3789.9167
$ws.Range("K136").Value = 6811.5
$ws.Range("L136").Value = 11369.7501
$ws.Range("M136").Value = -4261.5
$ws.Range("N136").Value = -16469.7501
$ws.Range("H139").Value = 51607.5
$ws.Range("J139").Value = 51607.5
$ws.Range("L139").Value = 51607.5
$ws.Range("N139").Value = -61887.5

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 506255.5
$ws.Range("I3").Value = 506255.5
$ws.Range("K3").Value = 506255.5
$ws.Range("M3").Value = -506141.5
$ws.Range("H5").Value = 9381.666999999999
$ws.Range("I5").Value = 25150
$ws.Range("J5").Value = 1497.5
$ws.Range("K5").Value = 25150
$ws.Range("L5").Value = 1497.5
$ws.Range("M5").Value = -25037
$ws.Range("N5").Value = -1723.5
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H107").Value = 1095.4667
$ws.Range("I107").Value = 1145.1428
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 1145.1428
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 774.8571999999999
$ws.Range("N107").Value = -4240

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4601.0435
$ws.Range("I31").Value = 2747.1292
$ws.Range("J31").Value = 6113.4473
$ws.Range("K31").Value = 2747.1292
$ws.Range("L31").Value = 6113.4473
$ws.Range("M31").Value = -2452.1292
$ws.Range("N31").Value = -6703.4473
$ws.Range("H34").Value = 4601.0435
$ws.Range("I34").Value = 2747.1292
$ws.Range("J34").Value = 6113.4473
$ws.Range("K34").Value = 2747.1292
$ws.Range("L34").Value = 6113.4473
$ws.Range("M34").Value = -2545.1292
$ws.Range("N34").Value = -6517.4473
$ws.Range("H37").Value = 8870
$ws.Range("J37").Value = 4300
$ws.Range("L37").Value = 4300
$ws.Range("N37").Value = -4514
$ws.Range("H58").Value = 20001600
$ws.Range("I58").Value = 29412960
$ws.Range("J58").Value = 2463.875
$ws.Range("K58").Value = 29412960
$ws.Range("L58").Value = 2463.875
$ws.Range("M58").Value = -29412757
$ws.Range("N58").Value = -2869.875
$ws.Range("H132").Value = 4764111.5
$ws.Range("I132").Value = 6946059
$ws.Range("J132").Value = 3499.818
$ws.Range("K132").Value = 20838177
$ws.Range("L132").Value = 10499.454
$ws.Range("M132").Value = -20835647
$ws.Range("N132").Value = -15559.454
$ws.Range("H136").Value = 20001600
$ws.Range("I136").Value = 29412960
$ws.Range("J136").Value = 2463.875
$ws.Range("K136").Value = 88238880
$ws.Range("L136").Value = 7391.625
$ws.Range("M136").Value = -88236330
$ws.Range("N136").Value = -12491.625

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2697.75
$ws.Range("I5").Value = 1583.4
$ws.Range("J5").Value = 4555
$ws.Range("K5").Value = 4750.200000000001
$ws.Range("L5").Value = 13665
$ws.Range("M5").Value = -4638.200000000001
$ws.Range("N5").Value = -13889
$ws.Range("H129").Value = 1015.3461
$ws.Range("J129").Value = 1809.4166
$ws.Range("L129").Value = 5428.2498
$ws.Range("N129").Value = -15428.2498
$ws.Range("H135").Value = 2697.75
$ws.Range("I135").Value = 1583.4
$ws.Range("J135").Value = 4555
$ws.Range("K135").Value = 14250.6
$ws.Range("L135").Value = 40995
$ws.Range("M135").Value = -11715.6
$ws.Range("N135").Value = -46065

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3681.3635
$ws.Range("I7").Value = 6990
$ws.Range("J7").Value = 3523.8096
$ws.Range("K7").Value = 6990
$ws.Range("L7").Value = 3523.8096
$ws.Range("M7").Value = -6878
$ws.Range("N7").Value = -3747.8096
$ws.Range("H61").Value = 1430.2858
$ws.Range("I61").Value = 891.55554
$ws.Range("K61").Value = 891.55554
$ws.Range("M61").Value = -689.55554
$ws.Range("H93").Value = 1941.5714
$ws.Range("I93").Value = 940.6
$ws.Range("K93").Value = 940.6
$ws.Range("M93").Value = 307.4
$ws.Range("H113").Value = 1430.2858
$ws.Range("I113").Value = 891.55554
$ws.Range("K113").Value = 891.55554
$ws.Range("M113").Value = 1278.44446
$ws.Range("H126").Value = 3681.3635
$ws.Range("I126").Value = 6990
$ws.Range("J126").Value = 3523.8096
$ws.Range("K126").Value = 20970
$ws.Range("L126").Value = 10571.4288
$ws.Range("M126").Value = -18500
$ws.Range("N126").Value = -15511.4288
$ws.Range("H132").Value = 6345.923
$ws.Range("I132").Value = 4118.8335
$ws.Range("J132").Value = 8254.857
$ws.Range("K132").Value = 12356.5005
$ws.Range("L132").Value = 24764.571
$ws.Range("M132").Value = -9826.500499999998
$ws.Range("N132").Value = -29824.571

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3064.724
$ws.Range("I136").Value = 811.5833
$ws.Range("J136").Value = 13879.8
$ws.Range("K136").Value = 2434.7499
$ws.Range("L136").Value = 41639.39999999999
$ws.Range("M136").Value = 115.2501000000002
$ws.Range("N136").Value = -46739.39999999999
